$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name labels (column A) per reordered list ---
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 25 de Abril de 2020 a las 19:52'
$ws.Cells.Item(36, 1).Value = 'Emiratos Arabes Unidos'
$ws.Cells.Item(37, 1).Value = 'Bielorrusia'
$ws.Cells.Item(38, 1).Value = 'Catar'
$ws.Cells.Item(62, 1).Value = 'Kazajistan'
$ws.Cells.Item(63, 1).Value = 'Barein'
$ws.Cells.Item(137, 1).Value = 'Maldivas'
$ws.Cells.Item(138, 1).Value = 'Brunei'
$ws.Cells.Item(164, 1).Value = 'Suazilandia'
$ws.Cells.Item(165, 1).Value = 'Benin'
$ws.Cells.Item(166, 1).Value = 'Guinea-Bisau'
$ws.Cells.Item(167, 1).Value = 'Nepal'
$ws.Cells.Item(168, 1).Value = 'Republica del Chad'
$ws.Cells.Item(169, 1).Value = 'Macao'
$ws.Cells.Item(170, 1).Value = 'Siria'
$ws.Cells.Item(211, 1).Value = 'Bonaire, San Eustaquio y Saba'
$ws.Cells.Item(212, 1).Value = 'Sudan del Sur'

# --- Update numeric statistics (columns B-H) ---
# Row 4
$ws.Cells.Item(4, 2).Value = 945249
$ws.Cells.Item(4, 3).Value = 20017
$ws.Cells.Item(4, 4).Value = 110834
$ws.Cells.Item(4, 5).Value = 781172
$ws.Cells.Item(4, 7).Value = 1050
$ws.Cells.Item(4, 8).Value = 53243
# Row 7
$ws.Cells.Item(7, 2).Value = 161488
$ws.Cells.Item(7, 3).Value = 1660
$ws.Cells.Item(7, 4).Value = 44594
$ws.Cells.Item(7, 5).Value = 94280
$ws.Cells.Item(7, 6).Value = 4725
$ws.Cells.Item(7, 7).Value = 369
$ws.Cells.Item(7, 8).Value = 22614
# Row 36
$ws.Cells.Item(36, 2).Value = 9813
$ws.Cells.Item(36, 3).Value = 532
$ws.Cells.Item(36, 4).Value = 1760
$ws.Cells.Item(36, 5).Value = 7989
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 64
# Row 37
$ws.Cells.Item(37, 2).Value = 9590
$ws.Cells.Item(37, 3).Value = 817
$ws.Cells.Item(37, 4).Value = 1573
$ws.Cells.Item(37, 5).Value = 7950
$ws.Cells.Item(37, 6).Value = 92
$ws.Cells.Item(37, 7).Value = 4
$ws.Cells.Item(37, 8).Value = 67
# Row 38
$ws.Cells.Item(38, 2).Value = 9358
$ws.Cells.Item(38, 3).Value = 833
$ws.Cells.Item(38, 4).Value = 929
$ws.Cells.Item(38, 5).Value = 8419
$ws.Cells.Item(38, 6).Value = 72
$ws.Cells.Item(38, 8).Value = 10
# Row 47
$ws.Cells.Item(47, 2).Value = 5926
$ws.Cells.Item(47, 3).Value = 177
$ws.Cells.Item(47, 4).Value = 822
$ws.Cells.Item(47, 5).Value = 4831
$ws.Cells.Item(47, 6).Value = 136
# Row 62
$ws.Cells.Item(62, 2).Value = 2601
$ws.Cells.Item(62, 3).Value = 185
$ws.Cells.Item(62, 4).Value = 646
$ws.Cells.Item(62, 5).Value = 1930
$ws.Cells.Item(62, 6).Value = 31
$ws.Cells.Item(62, 8).Value = 25
# Row 63
$ws.Cells.Item(63, 2).Value = 2588
$ws.Cells.Item(63, 3).Value = 70
$ws.Cells.Item(63, 4).Value = 1160
$ws.Cells.Item(63, 5).Value = 1420
$ws.Cells.Item(63, 6).Value = 2
$ws.Cells.Item(63, 8).Value = 8
# Row 137
$ws.Cells.Item(137, 2).Value = 141
$ws.Cells.Item(137, 3).Value = 12
$ws.Cells.Item(137, 4).Value = 17
$ws.Cells.Item(137, 5).Value = 124
$ws.Cells.Item(137, 8).Value = 0
# Row 138
$ws.Cells.Item(138, 2).Value = 138
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 121
$ws.Cells.Item(138, 5).Value = 16
$ws.Cells.Item(138, 8).Value = 1
# Row 164
$ws.Cells.Item(164, 2).Value = 56
$ws.Cells.Item(164, 3).Value = 20
$ws.Cells.Item(164, 4).Value = 10
$ws.Cells.Item(164, 5).Value = 45
# Row 165
$ws.Cells.Item(165, 2).Value = 54
$ws.Cells.Item(165, 4).Value = 27
$ws.Cells.Item(165, 5).Value = 26
$ws.Cells.Item(165, 8).Value = 1
# Row 166
$ws.Cells.Item(166, 2).Value = 52
$ws.Cells.Item(166, 4).Value = 3
$ws.Cells.Item(166, 5).Value = 49
# Row 167
$ws.Cells.Item(167, 2).Value = 49
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 12
$ws.Cells.Item(167, 5).Value = 37
# Row 168
$ws.Cells.Item(168, 2).Value = 46
$ws.Cells.Item(168, 3).Value = 6
$ws.Cells.Item(168, 4).Value = 15
$ws.Cells.Item(168, 5).Value = 31
$ws.Cells.Item(168, 6).Value = 0
# Row 169
$ws.Cells.Item(169, 2).Value = 45
$ws.Cells.Item(169, 4).Value = 27
$ws.Cells.Item(169, 5).Value = 18
$ws.Cells.Item(169, 6).Value = 1
$ws.Cells.Item(169, 8).Value = 0
# Row 170
$ws.Cells.Item(170, 2).Value = 42
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 6
$ws.Cells.Item(170, 5).Value = 33
$ws.Cells.Item(170, 8).Value = 3
# Row 172
$ws.Cells.Item(172, 4).Value = 13
